# Insert a new data row at row 141 (Fruta / hortaliza, semanal update),
# shifting the existing rows 141-245 down to 142-246.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(141).Insert()

# Copy the date cell's number format from the row below (now row 142,
# formerly row 141) so the new row's date cell is styled consistently.
$ws.Range("D141").NumberFormat = $ws.Range("D142").NumberFormat

$ws.Range("A141").Value = 5
$ws.Range("B141").Value = "Macroferia Regional de Talca"
$ws.Range("C141").Value = "Maule"
$ws.Range("D141").Value = 44574
$ws.Range("E141").Value = 7
$ws.Range("F141").Value = "Fruta"
$ws.Range("G141").Value = 100102
$ws.Range("H141").Value = "Cítricos"
$ws.Range("I141").Value = 100102004
$ws.Range("J141").Value = "Mandarina"
$ws.Range("K141").Value = "Murcott"
$ws.Range("L141").Value = "Primera"
$ws.Range("M141").Value = 350
$ws.Range("N141").Value = 7000
$ws.Range("O141").Value = 7000
$ws.Range("P141").Value = 7000
$ws.Range("Q141").Value = "`$/caja 18 kilos"
$ws.Range("R141").Value = "Región de O'Higgins"
$ws.Range("S141").Value = 389
$ws.Range("T141").Value = 18
